$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 85.388885
$ws.Range("I11").Value = 85.388885
$ws.Range("K11").Value = 85.388885
$ws.Range("M11").Value = 54.611115

$ws.Range("H12").Value = 374
$ws.Range("J12").Value = 415.33334
$ws.Range("L12").Value = 415.33334
$ws.Range("N12").Value = -755.33334

$ws.Range("H18").Value = 9953.308000000001
$ws.Range("I18").Value = 1349.75
$ws.Range("K18").Value = 1349.75
$ws.Range("M18").Value = -1065.75

$ws.Range("H28").Value = 716.6842
$ws.Range("I28").Value = 389.2353
$ws.Range("K28").Value = 389.2353
$ws.Range("M28").Value = 95.7647

$ws.Range("H33").Value = 209.6207
$ws.Range("I33").Value = 228.10527
$ws.Range("J33").Value = 174.5
$ws.Range("K33").Value = 228.10527
$ws.Range("L33").Value = 174.5
$ws.Range("M33").Value = 0.8947300000000098
$ws.Range("N33").Value = -632.5

$ws.Range("H74").Value = 7701.778
$ws.Range("I74").Value = 4257.25
$ws.Range("K74").Value = 4257.25
$ws.Range("M74").Value = -3321.25

$ws.Range("H77").Value = 7701.778
$ws.Range("I77").Value = 4257.25
$ws.Range("K77").Value = 21286.25
$ws.Range("M77").Value = -16606.25

$ws.Range("H92").Value = 3012.611
$ws.Range("I92").Value = 3391.9167
$ws.Range("K92").Value = 3391.9167
$ws.Range("M92").Value = -2143.9167

$ws.Range("H127").Value = 3859
$ws.Range("I127").Value = 3859
$ws.Range("K127").Value = 11577
$ws.Range("M127").Value = -6617

$ws.Range("H137").Value = 1103062.9
$ws.Range("I137").Value = 1544289.1
$ws.Range("J137").Value = 989604.75
$ws.Range("K137").Value = 4632867.300000001
$ws.Range("L137").Value = 2968814.25
$ws.Range("M137").Value = -4630317.300000001
$ws.Range("N137").Value = -2973914.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 3061.75
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 3061.75
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -3889.75

$ws.Range("H132").Value = 458502.03
$ws.Range("I132").Value = 518413.44
$ws.Range("J132").Value = 9166.333000000001
$ws.Range("K132").Value = 1555240.32
$ws.Range("L132").Value = 27498.999
$ws.Range("M132").Value = -1552710.32
$ws.Range("N132").Value = -32558.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 5413.5
$ws.Range("I80").Value = 50000
$ws.Range("J80").Value = 1983.7693
$ws.Range("K80").Value = 50000
$ws.Range("L80").Value = 1983.7693
$ws.Range("M80").Value = -49002
$ws.Range("N80").Value = -3979.7693

$ws.Range("H83").Value = 5413.5
$ws.Range("I83").Value = 50000
$ws.Range("J83").Value = 1983.7693
$ws.Range("K83").Value = 250000
$ws.Range("L83").Value = 9918.8465
$ws.Range("M83").Value = -245008
$ws.Range("N83").Value = -19902.8465

$ws.Range("H134").Value = 430871.66
$ws.Range("I134").Value = 565174.25
$ws.Range("K134").Value = 1695522.75
$ws.Range("M134").Value = -1692987.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1015.25
$ws.Range("I22").Value = 929.5
$ws.Range("K22").Value = 929.5
$ws.Range("M22").Value = -579.5

$ws.Range("H39").Value = 3833.3333
$ws.Range("I39").Value = 3833.3333
$ws.Range("K39").Value = 3833.3333
$ws.Range("M39").Value = -3442.3333

$ws.Range("H49").Value = 3833.3333
$ws.Range("I49").Value = 3833.3333
$ws.Range("K49").Value = 3833.3333
$ws.Range("M49").Value = -3651.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 178600.36
$ws.Range("I2").Value = 227294.86
$ws.Range("J2").Value = 53.833332
$ws.Range("K2").Value = 1363769.16
$ws.Range("L2").Value = 322.999992
$ws.Range("M2").Value = -1363656.16
$ws.Range("N2").Value = -548.999992

$ws.Range("H68").Value = 8619.846
$ws.Range("I68").Value = 2173.8333
$ws.Range("K68").Value = 6521.499899999999
$ws.Range("M68").Value = -5710.499899999999

$ws.Range("H71").Value = 8619.846
$ws.Range("I71").Value = 2173.8333
$ws.Range("K71").Value = 19564.4997
$ws.Range("M71").Value = -15508.4997

$ws.Range("H131").Value = 18680.615
$ws.Range("I131").Value = 849
$ws.Range("J131").Value = 26605.777
$ws.Range("K131").Value = 2547
$ws.Range("L131").Value = 79817.33099999999
$ws.Range("M131").Value = 2493
$ws.Range("N131").Value = -89897.33099999999

$ws.Range("H137").Value = 5395.0625
$ws.Range("J137").Value = 7006.5
$ws.Range("L137").Value = 21019.5
$ws.Range("N137").Value = -31219.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 14310000
$ws.Range("J24").Value = 28333.334
$ws.Range("L24").Value = 28333.334
$ws.Range("N24").Value = -28679.334

$ws.Range("H126").Value = 1045035.06
$ws.Range("J126").Value = 3944.889
$ws.Range("L126").Value = 11834.667
$ws.Range("N126").Value = -16774.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 86742.5
$ws.Range("I22").Value = 169491.33
$ws.Range("K22").Value = 169491.33
$ws.Range("M22").Value = -169196.33

$ws.Range("H27").Value = 86742.5
$ws.Range("I27").Value = 169491.33
$ws.Range("K27").Value = 169491.33
$ws.Range("M27").Value = -169384.33

$ws.Range("H74").Value = 62155.668

$ws.Range("H77").Value = 62155.668

$ws.Range("H122").Value = 3399.9429
$ws.Range("I122").Value = 3153.074
$ws.Range("K122").Value = 9459.222
$ws.Range("M122").Value = -7009.222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 7070
$ws.Range("I51").Value = 7070
$ws.Range("K51").Value = 7070
$ws.Range("M51").Value = -6560

$ws.Range("H54").Value = 44767.25
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 44767.25
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -45807.25

$ws.Range("H81").Value = 1462.4546
$ws.Range("I81").Value = 1365.8334
$ws.Range("J81").Value = 1578.4
$ws.Range("K81").Value = 2731.6668
$ws.Range("L81").Value = 3156.8
$ws.Range("M81").Value = -1670.6668
$ws.Range("N81").Value = -5278.8

$ws.Range("H84").Value = 1462.4546
$ws.Range("I84").Value = 1365.8334
$ws.Range("J84").Value = 1578.4
$ws.Range("K84").Value = 13658.334
$ws.Range("L84").Value = 15784
$ws.Range("M84").Value = -8354.333999999999
$ws.Range("N84").Value = -26392
